# Updates cryptos list: price and volume(1h) refresh, plus a few rank-order swaps
# (rows 7/8 swap XRP/USDC order; rows 42-47 reorder several tokens)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='67.284.04'},
    @{Cell='E2'; Value='  +0.50%  '},
    @{Cell='D3'; Value='3.476.25'},
    @{Cell='E3'; Value='  -0.68%  '},
    @{Cell='E4'; Value='  +0.00%  '},
    @{Cell='D5'; Value='594.86'},
    @{Cell='E5'; Value='  +0.05%  '},
    @{Cell='D6'; Value='179.90'},
    @{Cell='E6'; Value='  +4.08%  '},
    @{Cell='B7'; Value='XRP'},
    @{Cell='C7'; Value='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'},
    @{Cell='D7'; Value='0.610'},
    @{Cell='E7'; Value='  +5.20%  '},
    @{Cell='B8'; Value='USDC'},
    @{Cell='C8'; Value='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'},
    @{Cell='D8'; Value='0.999'},
    @{Cell='E8'; Value='  +0.02%  '},
    @{Cell='D9'; Value='3.475.31'},
    @{Cell='E9'; Value='  -0.70%  '},
    @{Cell='E10'; Value='  +4.86%  '},
    @{Cell='E11'; Value='  -2.36%  '},
    @{Cell='D12'; Value='0.433'},
    @{Cell='E12'; Value='  +0.80%  '},
    @{Cell='D13'; Value='4.073.41'},
    @{Cell='E13'; Value='  -0.48%  '},
    @{Cell='D14'; Value='31.84'},
    @{Cell='E14'; Value='  +8.10%  '},
    @{Cell='E15'; Value='  +0.04%  '},
    @{Cell='D16'; Value='67.269.31'},
    @{Cell='E16'; Value='  +0.54%  '},
    @{Cell='D17'; Value='0.0000177'},
    @{Cell='E17'; Value='  -0.75%  '},
    @{Cell='D18'; Value='3.474.85'},
    @{Cell='E18'; Value='  -0.57%  '},
    @{Cell='D19'; Value='6.27'},
    @{Cell='E19'; Value='  +0.19%  '},
    @{Cell='D20'; Value='14.19'},
    @{Cell='E20'; Value='  -0.79%  '},
    @{Cell='D21'; Value='389.33'},
    @{Cell='E21'; Value='  -0.33%  '},
    @{Cell='D22'; Value='7.96'},
    @{Cell='E22'; Value='  +0.41%  '},
    @{Cell='D23'; Value='72.73'},
    @{Cell='E23'; Value='  -0.66%  '},
    @{Cell='D24'; Value='0.999'},
    @{Cell='E24'; Value='  -0.14%  '},
    @{Cell='E25'; Value='  +1.30%  '},
    @{Cell='E26'; Value='  +0.58%  '},
    @{Cell='D27'; Value='0.0000123'},
    @{Cell='E27'; Value='  +1.43%  '},
    @{Cell='D28'; Value='10.27'},
    @{Cell='E28'; Value='  +1.56%  '},
    @{Cell='E29'; Value='  -2.67%  '},
    @{Cell='E30'; Value='  +0.47%  '},
    @{Cell='D31'; Value='6.21'},
    @{Cell='E31'; Value='  +1.69%  '},
    @{Cell='D32'; Value='1.42'},
    @{Cell='E32'; Value='  -0.20%  '},
    @{Cell='E33'; Value='  +0.29%  '},
    @{Cell='D34'; Value='23.51'},
    @{Cell='E34'; Value='  -0.38%  '},
    @{Cell='D35'; Value='7.37'},
    @{Cell='E35'; Value='  +0.42%  '},
    @{Cell='E36'; Value='  -0.01%  '},
    @{Cell='E37'; Value='  -0.08%  '},
    @{Cell='D38'; Value='161.79'},
    @{Cell='E38'; Value='  -0.85%  '},
    @{Cell='D39'; Value='0.884'},
    @{Cell='E39'; Value='  +0.94%  '},
    @{Cell='D41'; Value='1.87'},
    @{Cell='E41'; Value='  -1.80%  '},
    @{Cell='B42'; Value='Filecoin'},
    @{Cell='C42'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{Cell='D42'; Value='4.70'},
    @{Cell='E42'; Value='  +1.63%  '},
    @{Cell='B43'; Value='RenderToken'},
    @{Cell='C43'; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'},
    @{Cell='D43'; Value='6.79'},
    @{Cell='E43'; Value='  -0.31%  '},
    @{Cell='B44'; Value='EnergySwap'},
    @{Cell='C44'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell='D44'; Value='26.13'},
    @{Cell='E44'; Value='  +0.74%  '},
    @{Cell='B45'; Value='Maker'},
    @{Cell='C45'; Value='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{Cell='D45'; Value='2.805.87'},
    @{Cell='E45'; Value='  -0.46%  '},
    @{Cell='B46'; Value='Hedera'},
    @{Cell='C46'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Cell='D46'; Value='0.0722'},
    @{Cell='E46'; Value='  -0.71%  '},
    @{Cell='B47'; Value='InjectiveProtocol'},
    @{Cell='C47'; Value='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'},
    @{Cell='D47'; Value='26.48'},
    @{Cell='E47'; Value='  -1.67%  '},
    @{Cell='D48'; Value='41.08'},
    @{Cell='E48'; Value='  -3.11%  '},
    @{Cell='E49'; Value='  -0.37%  '},
    @{Cell='D50'; Value='331.20'},
    @{Cell='E50'; Value='  -2.51%  '},
    @{Cell='E51'; Value='  -2.10%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "0.610", "67.284.04")
    # are not coerced into doubles and keep exact text (trailing zeros, dots-as-
    # thousands-separators). Reset the style back to Normal afterward so no stray
    # cell formatting is introduced beyond the text value itself.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
